$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Unmerge the old A147:A149 "joint segmentation and detection" heading block - in the
#    new layout column A is no longer a merged heading there.
$ws.Range("A147:A149").UnMerge()

# 2. Insert a new blank row at 148 (pushes old rows 148-152 down to 149-153). The
#    A150:A152 "Lane detection" merge automatically becomes A151:A153.
$ws.Rows.Item(148).Insert()

# 3. Row 147 keeps "joint segmentation and detection" in A147 (now a plain, non-centered
#    wrap-text cell) and gets a new paper title in C147; grow the row height.
$ws.Range("A147").Value = "joint segmentation and detection"
$ws.Range("A147").HorizontalAlignment = 1
$ws.Range("C147").Value = "Multi-Task Learning Using Uncertainty to Weigh Losses for Scene Geometry and Semantics"
$ws.Range("D147").Value = ""
$ws.Rows.Item(147).RowHeight = 56

# 4. New row 148 holds the note for "An overview of multi-task learning in deep neural
#    networks" plus its comment in column D. Clear column A entirely (no heading there).
$ws.Range("A148").ClearFormats()
$ws.Range("B148").Value = 5
$ws.Range("C148").Value = "An overview of multi-task learning in deep neural networks"
$ws.Range("D148").Value = "why does MTL work: 1. implicit data augmentation; 2. attention focusing; 3. eavesdropping; 4. representation bias; 5. regularization`n到底如何共享backbone是需要尝试的"
$ws.Rows.Item(148).RowHeight = 42

# Style B148 like the other "priority" cells inside the bordered table (fill + border),
# and C148/D148 like the borderless wrap-text notes used below the table.
$ws.Range("B148").Interior.Color = $ws.Range("B146").Interior.Color
$ws.Range("B148").Borders.LineStyle = $ws.Range("B146").Borders.Item(1).LineStyle
$ws.Range("B148").HorizontalAlignment = -4108
$ws.Range("B148").VerticalAlignment = -4108
$ws.Range("C148:D148").WrapText = $true
$ws.Range("C148:D148").VerticalAlignment = -4108
$ws.Range("C148:D148").HorizontalAlignment = 1

# 5. Rows 149-150 (old 148-149 content) are no longer part of a merged/centered column A
#    heading - switch them to plain, non-centered wrap-text cells.
$ws.Range("A149").HorizontalAlignment = 1
$ws.Range("A150").HorizontalAlignment = 1
